# Update the Stock column (E) for the existing medicine records.
# pk_MedId 16 (row 3): Stock 12 -> 3
# pk_MedId 17 (row 4): Stock 10 -> 4
# pk_MedId 18 (row 5): Stock 10 -> 0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# Forcing NumberFormat to text before assignment keeps the value stored as
# a shared string (matching the exporter's original "numeric-looking text"
# convention) instead of Excel's default numeric auto-coercion; clearing
# the format afterwards drops the temporary text format so the cell keeps
# its original (default) style.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3"
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4"
$ws.Range("E4").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0"
$ws.Range("E5").ClearFormats()
